# "updated with free crm"
# The two existing test cases (LoginTest / ValidateCRMTest, which drove a
# CRM login flow) are replaced with two "search Google" demo test cases,
# and their TestData blocks are renamed/duplicated accordingly.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("TestCases")
$ws2 = $wb.Worksheets.Item("TestData")

# --- Sheet "TestCases": rename the two test case identifiers ---
$ws1.Range("A2").Value = "searchGoogle"
$ws1.Range("A3").Value = "searchGoogle2"

# --- Sheet "TestData": rename the first data block (for searchGoogle) ---
$ws2.Range("A1").Value = "searchGoogle"
$ws2.Range("B2").Value = "SearchKeyword"
$ws2.Range("B3").Value = "Selenium Automation demo 1"

# --- Sheet "TestData": append a second data block (for searchGoogle2) ---
$ws2.Range("A6").Value = "searchGoogle2"
$ws2.Range("A7").Value = "Runmode"
$ws2.Range("B7").Value = "SearchKeyword"
$ws2.Range("C7").Value = "password"
$ws2.Range("D7").Value = "browser"
$ws2.Range("A8").Value = "Y"
$ws2.Range("B8").Value = "Selenium Automation demo 2"
$ws2.Range("C8").Value = "Admin@123"
$ws2.Range("D8").Value = "chrome"

# --- Selection / active-sheet state: TestData shows B2 selected; the
#     workbook re-opens on TestCases with A4 selected and tab-active ---
$ws2.Range("B2").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("A4").Select() | Out-Null
